$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B ("Date") before the existing WeekNumber column,
# shifting old B:E to C:F.
$ws.Columns("B").Insert()

# New column B must hold its values as TEXT (shared strings), matching the
# original author's data (dates typed as plain text, not date serials).
# Setting the Text number format before assigning values keeps them as
# strings instead of Excel auto-converting "2015-01-28" style input into a
# date serial number.
$ws.Columns("B").NumberFormat = "@"

$ws.Range("B1").Value = "Date"

$ws.Range("B2").Value = "2015-01-28"
$ws.Range("B3").Value = "2015-01-28"
$ws.Range("B4").Value = "2015-01-28"
$ws.Range("B5").Value = "2015-01-28"
$ws.Range("B6").Value = "2015-01-28"
$ws.Range("B7").Value = "2015-02-05"
$ws.Range("B8").Value = "2015-02-05"
$ws.Range("B9").Value = "2015-02-05"
$ws.Range("B10").Value = "2015-01-29"
$ws.Range("B11").Value = "2015-02-04"
$ws.Range("B12").Value = "2015-02-11"
$ws.Range("B13").Value = "2015-02-11"
$ws.Range("B14").Value = "2015-02-05"
$ws.Range("B15").Value = "2015-02-18"
$ws.Range("B16").Value = "2015-02-25"
$ws.Range("B17").Value = "2015-02-25"
$ws.Range("B18").Value = "2015-03-04"
$ws.Range("B19").Value = "2015-03-11"
$ws.Range("B20").Value = "2015-03-11"
$ws.Range("B21").Value = "2015-03-18"
$ws.Range("B22").Value = "2015-03-18"
$ws.Range("B23").Value = "2015-03-25"
$ws.Range("B24").Value = "2015-04-01"
$ws.Range("B25").Value = "2015-04-01"
$ws.Range("B26").Value = "2015-04-01"
$ws.Range("B27").Value = "2015-04-08"
$ws.Range("B28").Value = "2015-04-08"
$ws.Range("B29").Value = "2015-04-08"
$ws.Range("B30").Value = "2015-04-15"
$ws.Range("B31").Value = "2015-04-15"
$ws.Range("B32").Value = "2015-04-15"
$ws.Range("B33").Value = "2015-04-15"
$ws.Range("B34").Value = "2015-04-22"
$ws.Range("B35").Value = "2015-04-22"
$ws.Range("B36").Value = "2015-04-22"
$ws.Range("B37").Value = "2015-04-29"

# Match the column width used by column A so the new Date column renders the
# same way.
$ws.Columns("B").ColumnWidth = $ws.Columns("A").ColumnWidth

# Match the author's final selection / active cell.
$ws.Range("B9").Select() | Out-Null
